$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill KH1:KM102 with the updated dataset (round 4 -> round 9 columns added,
# KH re-stated with its existing value so the whole block can be written as one array).
$arr = New-Object 'object[,]' 102,6
$arr[0,0] = 11002
$arr[0,1] = 11007
$arr[0,2] = 11016
$arr[0,3] = 11023
$arr[0,4] = 11032
$arr[0,5] = 11043
$arr[1,0] = 2023
$arr[1,1] = 2023
$arr[1,2] = 2023
$arr[1,3] = 2023
$arr[1,4] = 2023
$arr[1,5] = 2023
$arr[2,0] = 4
$arr[2,1] = 5
$arr[2,2] = 6
$arr[2,3] = 7
$arr[2,4] = 8
$arr[2,5] = 9
$arr[3,0] = 1
$arr[3,1] = 0
$arr[3,2] = 1
$arr[3,3] = 0
$arr[3,4] = 1
$arr[3,5] = 1
$arr[4,0] = 1
$arr[4,1] = 0
$arr[4,2] = 0
$arr[4,3] = 0
$arr[4,4] = 0
$arr[4,5] = 0
$arr[5,0] = 117
$arr[5,1] = 80
$arr[5,2] = 98
$arr[5,3] = 113
$arr[5,4] = 69
$arr[5,5] = 62
$arr[6,0] = 89
$arr[6,1] = 79
$arr[6,2] = 117
$arr[6,3] = 59
$arr[6,4] = 98
$arr[6,5] = 82
$arr[7,0] = 28
$arr[7,1] = 1
$arr[7,2] = -19
$arr[7,3] = 54
$arr[7,4] = -29
$arr[7,5] = -20
$arr[8,0] = 1
$arr[8,1] = 1
$arr[8,2] = 0
$arr[8,3] = 1
$arr[8,4] = 0
$arr[8,5] = 0
$arr[9,0] = 8
$arr[9,1] = 15
$arr[9,2] = 3
$arr[9,3] = 2
$arr[9,4] = 16
$arr[9,5] = 5
$arr[10,0] = 206
$arr[10,1] = 210
$arr[10,2] = 217
$arr[10,3] = 237
$arr[10,4] = 199
$arr[10,5] = 186
$arr[11,0] = 151
$arr[11,1] = 173
$arr[11,2] = 133
$arr[11,3] = 151
$arr[11,4] = 142
$arr[11,5] = 164
$arr[12,0] = 357
$arr[12,1] = 383
$arr[12,2] = 350
$arr[12,3] = 388
$arr[12,4] = 341
$arr[12,5] = 350
$arr[13,0] = 1.36
$arr[13,1] = 1.21
$arr[13,2] = 1.63
$arr[13,3] = 1.57
$arr[13,4] = 1.4
$arr[13,5] = 1.13
$arr[14,0] = 83
$arr[14,1] = 95
$arr[14,2] = 106
$arr[14,3] = 113
$arr[14,4] = 70
$arr[14,5] = 74
$arr[15,0] = 47
$arr[15,1] = 59
$arr[15,2] = 68
$arr[15,3] = 56
$arr[15,4] = 72
$arr[15,5] = 69
$arr[16,0] = 15
$arr[16,1] = 33
$arr[16,2] = 28
$arr[16,3] = 35
$arr[16,4] = 27
$arr[16,5] = 30
$arr[17,0] = 19
$arr[17,1] = 9
$arr[17,2] = 12
$arr[17,3] = 16
$arr[17,4] = 19
$arr[17,5] = 10
$arr[18,0] = 16
$arr[18,1] = 20
$arr[18,2] = 19
$arr[18,3] = 21
$arr[18,4] = 27
$arr[18,5] = 12
$arr[19,0] = 18
$arr[19,1] = 11
$arr[19,2] = 15
$arr[19,3] = 17
$arr[19,4] = 10
$arr[19,5] = 9
$arr[20,0] = 15
$arr[20,1] = 7
$arr[20,2] = 10
$arr[20,3] = 13
$arr[20,4] = 8
$arr[20,5] = 4
$arr[21,0] = 7
$arr[21,1] = 11
$arr[21,2] = 4
$arr[21,3] = 7
$arr[21,4] = 7
$arr[21,5] = 6
$arr[22,0] = 2
$arr[22,1] = 3
$arr[22,2] = 4
$arr[22,3] = 4
$arr[22,4] = 2
$arr[22,5] = 2
$arr[23,0] = 27
$arr[23,1] = 25
$arr[23,2] = 23
$arr[23,3] = 28
$arr[23,4] = 19
$arr[23,5] = 17
$arr[24,0] = 66.7
$arr[24,1] = 44
$arr[24,2] = 65.2
$arr[24,3] = 60.7
$arr[24,4] = 52.6
$arr[24,5] = 52.9
$arr[25,0] = 19.83
$arr[25,1] = 34.82
$arr[25,2] = 23.33
$arr[25,3] = 22.82
$arr[25,4] = 34.1
$arr[25,5] = 38.89
$arr[26,0] = 13.22
$arr[26,1] = 15.32
$arr[26,2] = 15.22
$arr[26,3] = 13.86
$arr[26,4] = 17.95
$arr[26,5] = 20.59
$arr[27,0] = 37
$arr[27,1] = 31
$arr[27,2] = 32
$arr[27,3] = 38
$arr[27,4] = 31
$arr[27,5] = 40
$arr[28,0] = 66
$arr[28,1] = 44
$arr[28,2] = 54
$arr[28,3] = 53
$arr[28,4] = 73
$arr[28,5] = 52
$arr[29,0] = 40
$arr[29,1] = 40
$arr[29,2] = 37
$arr[29,3] = 40
$arr[29,4] = 52
$arr[29,5] = 45
$arr[30,0] = 69
$arr[30,1] = 48
$arr[30,2] = 51
$arr[30,3] = 53
$arr[30,4] = 42
$arr[30,5] = 40
$arr[31,0] = 2.56
$arr[31,1] = 1.92
$arr[31,2] = 2.22
$arr[31,3] = 1.89
$arr[31,4] = 2.21
$arr[31,5] = 2.35
$arr[32,0] = 3.83
$arr[32,1] = 4.36
$arr[32,2] = 3.4
$arr[32,3] = 3.12
$arr[32,4] = 4.2
$arr[32,5] = 4.44
$arr[33,0] = 36.2
$arr[33,1] = 45.8
$arr[33,2] = 37.3
$arr[33,3] = 45.3
$arr[33,4] = 40.5
$arr[33,5] = 37.5
$arr[34,0] = 26.1
$arr[34,1] = 22.9
$arr[34,2] = 29.4
$arr[34,3] = 32.1
$arr[34,4] = 23.8
$arr[34,5] = 22.5
$arr[35,0] = 188.5
$arr[35,1] = 188.5
$arr[35,2] = 188.6
$arr[35,3] = 188.3
$arr[35,4] = 188.2
$arr[35,5] = 188.4
$arr[36,0] = 87.3
$arr[36,1] = 87.3
$arr[36,2] = 87
$arr[36,3] = 86.09999999999999
$arr[36,4] = 86.2
$arr[36,5] = 86.59999999999999
$arr[37,0] = 25.8
$arr[37,1] = 25.8
$arr[37,2] = 25.16
$arr[37,3] = 24.74
$arr[37,4] = 25.16
$arr[37,5] = 25.66
$arr[38,0] = 87.59999999999999
$arr[38,1] = 88.59999999999999
$arr[38,2] = 89.40000000000001
$arr[38,3] = 80.8
$arr[38,4] = 93
$arr[38,5] = 108.3
$arr[39,0] = 9
$arr[39,1] = 9
$arr[39,2] = 10
$arr[39,3] = 10
$arr[39,4] = 9
$arr[39,5] = 7
$arr[40,0] = 7
$arr[40,1] = 7
$arr[40,2] = 6
$arr[40,3] = 7
$arr[40,4] = 7
$arr[40,5] = 6
$arr[41,0] = 2
$arr[41,1] = 2
$arr[41,2] = 2
$arr[41,3] = 2
$arr[41,4] = 2
$arr[41,5] = 3
$arr[42,0] = 5
$arr[42,1] = 5
$arr[42,2] = 5
$arr[42,3] = 4
$arr[42,4] = 5
$arr[42,5] = 7
$arr[43,0] = 135
$arr[43,1] = 121
$arr[43,2] = 117
$arr[43,3] = 132
$arr[43,4] = 152
$arr[43,5] = 136
$arr[44,0] = 225
$arr[44,1] = 259
$arr[44,2] = 217
$arr[44,3] = 248
$arr[44,4] = 182
$arr[44,5] = 211
$arr[45,0] = 264
$arr[45,1] = 290
$arr[45,2] = 264
$arr[45,3] = 294
$arr[45,4] = 217
$arr[45,5] = 244
$arr[46,0] = 73.90000000000001
$arr[46,1] = 75.7
$arr[46,2] = 75.40000000000001
$arr[46,3] = 75.8
$arr[46,4] = 63.6
$arr[46,5] = 69.7
$arr[47,0] = 66
$arr[47,1] = 44
$arr[47,2] = 54
$arr[47,3] = 53
$arr[47,4] = 73
$arr[47,5] = 52
$arr[48,0] = 12
$arr[48,1] = 5
$arr[48,2] = 9
$arr[48,3] = 8
$arr[48,4] = 10
$arr[48,5] = 11
$arr[49,0] = 14
$arr[49,1] = 15
$arr[49,2] = 13
$arr[49,3] = 15
$arr[49,4] = 13
$arr[49,5] = 10
$arr[50,0] = 37
$arr[50,1] = 31
$arr[50,2] = 32
$arr[50,3] = 38
$arr[50,4] = 31
$arr[50,5] = 40
$arr[51,0] = 40
$arr[51,1] = 40
$arr[51,2] = 37
$arr[51,3] = 40
$arr[51,4] = 52
$arr[51,5] = 45
$arr[52,0] = 41
$arr[52,1] = 51
$arr[52,2] = 42
$arr[52,3] = 40
$arr[52,4] = 68
$arr[52,5] = 46
$arr[53,0] = 5
$arr[53,1] = 6
$arr[53,2] = 6
$arr[53,3] = 17
$arr[53,4] = 6
$arr[53,5] = 5
$arr[54,0] = 15
$arr[54,1] = 7
$arr[54,2] = 10
$arr[54,3] = 13
$arr[54,4] = 8
$arr[54,5] = 4
$arr[55,0] = 83.3
$arr[55,1] = 63.6
$arr[55,2] = 66.7
$arr[55,3] = 76.5
$arr[55,4] = 80
$arr[55,5] = 44.4
$arr[56,0] = 218
$arr[56,1] = 226
$arr[56,2] = 189
$arr[56,3] = 195
$arr[56,4] = 217
$arr[56,5] = 221
$arr[57,0] = 134
$arr[57,1] = 130
$arr[57,2] = 192
$arr[57,3] = 106
$arr[57,4] = 131
$arr[57,5] = 156
$arr[58,0] = 352
$arr[58,1] = 356
$arr[58,2] = 381
$arr[58,3] = 301
$arr[58,4] = 348
$arr[58,5] = 377
$arr[59,0] = 1.63
$arr[59,1] = 1.74
$arr[59,2] = 0.98
$arr[59,3] = 1.84
$arr[59,4] = 1.66
$arr[59,5] = 1.42
$arr[60,0] = 92
$arr[60,1] = 106
$arr[60,2] = 70
$arr[60,3] = 90
$arr[60,4] = 61
$arr[60,5] = 113
$arr[61,0] = 37
$arr[61,1] = 67
$arr[61,2] = 52
$arr[61,3] = 52
$arr[61,4] = 73
$arr[61,5] = 62
$arr[62,0] = 42
$arr[62,1] = 30
$arr[62,2] = 41
$arr[62,3] = 30
$arr[62,4] = 47
$arr[62,5] = 55
$arr[63,0] = 16
$arr[63,1] = 20
$arr[63,2] = 19
$arr[63,3] = 21
$arr[63,4] = 28
$arr[63,5] = 12
$arr[64,0] = 19
$arr[64,1] = 9
$arr[64,2] = 12
$arr[64,3] = 16
$arr[64,4] = 19
$arr[64,5] = 10
$arr[65,0] = 13
$arr[65,1] = 12
$arr[65,2] = 17
$arr[65,3] = 8
$arr[65,4] = 14
$arr[65,5] = 12
$arr[66,0] = 10
$arr[66,1] = 7
$arr[66,2] = 12
$arr[66,3] = 7
$arr[66,4] = 6
$arr[66,5] = 6
$arr[67,0] = 10
$arr[67,1] = 5
$arr[67,2] = 13
$arr[67,3] = 9
$arr[67,4] = 11
$arr[67,5] = 9
$arr[68,0] = 1
$arr[68,1] = 2
$arr[68,2] = 2
$arr[68,3] = 2
$arr[68,4] = 3
$arr[68,5] = 1
$arr[69,0] = 24
$arr[69,1] = 19
$arr[69,2] = 32
$arr[69,3] = 19
$arr[69,4] = 28
$arr[69,5] = 22
$arr[70,0] = 54.2
$arr[70,1] = 63.2
$arr[70,2] = 53.1
$arr[70,3] = 42.1
$arr[70,4] = 50
$arr[70,5] = 54.5
$arr[71,0] = 27.08
$arr[71,1] = 29.67
$arr[71,2] = 22.41
$arr[71,3] = 37.62
$arr[71,4] = 24.86
$arr[71,5] = 31.42
$arr[72,0] = 14.67
$arr[72,1] = 18.74
$arr[72,2] = 11.91
$arr[72,3] = 15.84
$arr[72,4] = 12.43
$arr[72,5] = 17.14
$arr[73,0] = 32
$arr[73,1] = 43
$arr[73,2] = 48
$arr[73,3] = 34
$arr[73,4] = 43
$arr[73,5] = 39
$arr[74,0] = 59
$arr[74,1] = 53
$arr[74,2] = 48
$arr[74,3] = 49
$arr[74,4] = 68
$arr[74,5] = 46
$arr[75,0] = 51
$arr[75,1] = 35
$arr[75,2] = 35
$arr[75,3] = 34
$arr[75,4] = 32
$arr[75,5] = 30
$arr[76,0] = 55
$arr[76,1] = 52
$arr[76,2] = 54
$arr[76,3] = 48
$arr[76,4] = 67
$arr[76,5] = 57
$arr[77,0] = 2.29
$arr[77,1] = 2.74
$arr[77,2] = 1.69
$arr[77,3] = 2.53
$arr[77,4] = 2.39
$arr[77,5] = 2.59
$arr[78,0] = 4.23
$arr[78,1] = 4.33
$arr[78,2] = 3.18
$arr[78,3] = 6
$arr[78,4] = 4.79
$arr[78,5] = 4.75
$arr[79,0] = 41.8
$arr[79,1] = 32.7
$arr[79,2] = 55.6
$arr[79,3] = 35.4
$arr[79,4] = 37.3
$arr[79,5] = 36.8
$arr[80,0] = 23.6
$arr[80,1] = 23.1
$arr[80,2] = 31.5
$arr[80,3] = 16.7
$arr[80,4] = 20.9
$arr[80,5] = 21.1
$arr[81,0] = 190
$arr[81,1] = 188.4
$arr[81,2] = 189.2
$arr[81,3] = 188.3
$arr[81,4] = 187.3
$arr[81,5] = 188.6
$arr[82,0] = 83.8
$arr[82,1] = 86.8
$arr[82,2] = 88.2
$arr[82,3] = 86.5
$arr[82,4] = 87.09999999999999
$arr[82,5] = 87.3
$arr[83,0] = 24
$arr[83,1] = 24.91
$arr[83,2] = 25.58
$arr[83,3] = 26.8
$arr[83,4] = 26.16
$arr[83,5] = 26.33
$arr[84,0] = 74.59999999999999
$arr[84,1] = 93.2
$arr[84,2] = 83.2
$arr[84,3] = 125.8
$arr[84,4] = 118.7
$arr[84,5] = 118.3
$arr[85,0] = 10
$arr[85,1] = 8
$arr[85,2] = 7
$arr[85,3] = 5
$arr[85,4] = 2
$arr[85,5] = 4
$arr[86,0] = 6
$arr[86,1] = 4
$arr[86,2] = 5
$arr[86,3] = 2
$arr[86,4] = 9
$arr[86,5] = 7
$arr[87,0] = 4
$arr[87,1] = 6
$arr[87,2] = 7
$arr[87,3] = 6
$arr[87,4] = 4
$arr[87,5] = 6
$arr[88,0] = 3
$arr[88,1] = 5
$arr[88,2] = 4
$arr[88,3] = 10
$arr[88,4] = 8
$arr[88,5] = 6
$arr[89,0] = 126
$arr[89,1] = 132
$arr[89,2] = 151
$arr[89,3] = 116
$arr[89,4] = 168
$arr[89,5] = 129
$arr[90,0] = 216
$arr[90,1] = 213
$arr[90,2] = 231
$arr[90,3] = 174
$arr[90,4] = 180
$arr[90,5] = 234
$arr[91,0] = 254
$arr[91,1] = 258
$arr[91,2] = 285
$arr[91,3] = 210
$arr[91,4] = 226
$arr[91,5] = 271
$arr[92,0] = 72.2
$arr[92,1] = 72.5
$arr[92,2] = 74.8
$arr[92,3] = 69.8
$arr[92,4] = 64.90000000000001
$arr[92,5] = 71.90000000000001
$arr[93,0] = 59
$arr[93,1] = 53
$arr[93,2] = 48
$arr[93,3] = 49
$arr[93,4] = 68
$arr[93,5] = 46
$arr[94,0] = 6
$arr[94,1] = 6
$arr[94,2] = 9
$arr[94,3] = 8
$arr[94,4] = 10
$arr[94,5] = 8
$arr[95,0] = 8
$arr[95,1] = 14
$arr[95,2] = 17
$arr[95,3] = 7
$arr[95,4] = 5
$arr[95,5] = 15
$arr[96,0] = 32
$arr[96,1] = 43
$arr[96,2] = 48
$arr[96,3] = 34
$arr[96,4] = 43
$arr[96,5] = 39
$arr[97,0] = 51
$arr[97,1] = 35
$arr[97,2] = 35
$arr[97,3] = 34
$arr[97,4] = 32
$arr[97,5] = 30
$arr[98,0] = 41
$arr[98,1] = 34
$arr[98,2] = 41
$arr[98,3] = 39
$arr[98,4] = 53
$arr[98,5] = 31
$arr[99,0] = 2
$arr[99,1] = 6
$arr[99,2] = 2
$arr[99,3] = 1
$arr[99,4] = 9
$arr[99,5] = 9
$arr[100,0] = 10
$arr[100,1] = 7
$arr[100,2] = 12
$arr[100,3] = 7
$arr[100,4] = 6
$arr[100,5] = 6
$arr[101,0] = 76.90000000000001
$arr[101,1] = 58.3
$arr[101,2] = 70.59999999999999
$arr[101,3] = 87.5
$arr[101,4] = 42.9
$arr[101,5] = 50

$ws.Range("KH1:KM102").Value = $arr

# Match the original formatting pattern: every column except the newest (KM)
# carries the sheet's standard cell style; KM is left as the freshly-appended,
# unstyled "latest round" column - same as KH was before this edit.
$ws.Range("KH1:KL102").Style = "Normal"

Write-Host "GWS_stats updated through round 9 (KH:KM)"
